# chore: update Sheets via scheduled runner
# Refreshes cached market-board derived figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) per-row across the eight crafting-class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 952511.56
$ws.Range("I19").Value = 1333418.2
$ws.Range("J19").Value = 244.75
$ws.Range("K19").Value = 1333418.2
$ws.Range("L19").Value = 244.75
$ws.Range("M19").Value = -1333243.2
$ws.Range("N19").Value = -594.75
$ws.Range("H51").Value = 5880.0835
$ws.Range("I51").Value = 3333.6667
$ws.Range("J51").Value = 6728.8887
$ws.Range("K51").Value = 3333.6667
$ws.Range("L51").Value = 6728.8887
$ws.Range("M51").Value = -2849.6667
$ws.Range("N51").Value = -7696.8887
$ws.Range("H75").Value = 39800
$ws.Range("J75").Value = 39800
$ws.Range("L75").Value = 39800
$ws.Range("N75").Value = -41672
$ws.Range("H78").Value = 39800
$ws.Range("J78").Value = 39800
$ws.Range("L78").Value = 119400
$ws.Range("N78").Value = -128760
$ws.Range("H132").Value = 142228.45
$ws.Range("I132").Value = 224968.06
$ws.Range("J132").Value = 6836.364
$ws.Range("K132").Value = 674904.1799999999
$ws.Range("L132").Value = 20509.092
$ws.Range("M132").Value = -672374.1799999999
$ws.Range("N132").Value = -25569.092
$ws.Range("H137").Value = 4485.4165
$ws.Range("I137").Value = 3811.926
$ws.Range("J137").Value = 6505.8887
$ws.Range("K137").Value = 11435.778
$ws.Range("L137").Value = 19517.6661
$ws.Range("M137").Value = -8885.778
$ws.Range("N137").Value = -24617.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 482.14706
$ws.Range("I2").Value = 482.375
$ws.Range("K2").Value = 482.375
$ws.Range("M2").Value = -369.375
$ws.Range("H32").Value = 7521.0356
$ws.Range("I32").Value = 4242.519
$ws.Range("K32").Value = 4242.519
$ws.Range("M32").Value = -3955.519
$ws.Range("H45").Value = 1588.6923
$ws.Range("I45").Value = 946.1429
$ws.Range("J45").Value = 2338.3333
$ws.Range("K45").Value = 946.1429
$ws.Range("L45").Value = 2338.3333
$ws.Range("M45").Value = -569.1429
$ws.Range("N45").Value = -3092.3333
$ws.Range("H61").Value = 1800.9656
$ws.Range("I61").Value = 1196.1818
$ws.Range("K61").Value = 1196.1818
$ws.Range("M61").Value = -984.1818000000001
$ws.Range("H107").Value = 20114
$ws.Range("J107").Value = 20114
$ws.Range("L107").Value = 20114
$ws.Range("N107").Value = -27794
$ws.Range("H109").Value = 26050
$ws.Range("J109").Value = 26050
$ws.Range("L109").Value = 26050
$ws.Range("N109").Value = -28824
$ws.Range("H110").Value = 1085.125
$ws.Range("I110").Value = 995.25
$ws.Range("K110").Value = 995.25
$ws.Range("M110").Value = 1049.75
$ws.Range("H116").Value = 482.14706
$ws.Range("I116").Value = 482.375
$ws.Range("K116").Value = 482.375
$ws.Range("M116").Value = 1811.625
$ws.Range("H132").Value = 2280.5454
$ws.Range("I132").Value = 1605.0435
$ws.Range("K132").Value = 4815.1305
$ws.Range("M132").Value = -2285.1305
$ws.Range("H136").Value = 1800.9656
$ws.Range("I136").Value = 1196.1818
$ws.Range("K136").Value = 3588.5454
$ws.Range("M136").Value = -1038.5454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 482.14706
$ws.Range("I3").Value = 482.375
$ws.Range("K3").Value = 482.375
$ws.Range("M3").Value = -368.375
$ws.Range("H20").Value = 5449
$ws.Range("I20").Value = 2114.08
$ws.Range("J20").Value = 11404.214
$ws.Range("K20").Value = 2114.08
$ws.Range("L20").Value = 11404.214
$ws.Range("M20").Value = -1867.08
$ws.Range("N20").Value = -11898.214

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2442.5894
$ws.Range("I31").Value = 857.08826
$ws.Range("K31").Value = 857.08826
$ws.Range("M31").Value = -562.08826
$ws.Range("H34").Value = 2442.5894
$ws.Range("I34").Value = 857.08826
$ws.Range("K34").Value = 857.08826
$ws.Range("M34").Value = -655.08826
$ws.Range("H97").Value = 34210
$ws.Range("J97").Value = 34210
$ws.Range("L97").Value = 34210
$ws.Range("N97").Value = -36192
$ws.Range("H105").Value = 2879.8
$ws.Range("I105").Value = 2874.2856
$ws.Range("J105").Value = 2892.6667
$ws.Range("K105").Value = 2874.2856
$ws.Range("L105").Value = 2892.6667
$ws.Range("M105").Value = -1127.2856
$ws.Range("N105").Value = -6386.6667
$ws.Range("H109").Value = 25000
$ws.Range("J109").Value = 25000
$ws.Range("L109").Value = 25000
$ws.Range("N109").Value = -27080
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H132").Value = 4197.933
$ws.Range("I132").Value = 3656
$ws.Range("J132").Value = 4739.8667
$ws.Range("K132").Value = 10968
$ws.Range("L132").Value = 14219.6001
$ws.Range("M132").Value = -8438
$ws.Range("N132").Value = -19279.6001
$ws.Range("H134").Value = 6297.0835
$ws.Range("I134").Value = 7828.2
$ws.Range("J134").Value = 3745.2222
$ws.Range("K134").Value = 23484.6
$ws.Range("L134").Value = 11235.6666
$ws.Range("M134").Value = -20949.6
$ws.Range("N134").Value = -16305.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3999
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3999
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 3999
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -4539
$ws.Range("H73").Value = 3999
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3999
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 3999
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -5871
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 2595.2856
$ws.Range("I122").Value = 1994.0454
$ws.Range("K122").Value = 5982.1362
$ws.Range("M122").Value = -3532.1362
$ws.Range("H132").Value = 2448.4473
$ws.Range("I132").Value = 1440.72
$ws.Range("J132").Value = 4386.385
$ws.Range("K132").Value = 4322.16
$ws.Range("L132").Value = 13159.155
$ws.Range("M132").Value = -1792.16
$ws.Range("N132").Value = -18219.155
$ws.Range("H136").Value = 12545.742
$ws.Range("J136").Value = 12944.267
$ws.Range("L136").Value = 38832.801
$ws.Range("N136").Value = -43932.801

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 46800
$ws.Range("J75").Value = 46800
$ws.Range("L75").Value = 46800
$ws.Range("N75").Value = -48672
$ws.Range("H78").Value = 46800
$ws.Range("J78").Value = 46800
$ws.Range("L78").Value = 140400
$ws.Range("N78").Value = -149760
$ws.Range("H122").Value = 3186.9556
$ws.Range("I122").Value = 2897.2727
$ws.Range("J122").Value = 3983.5833
$ws.Range("K122").Value = 8691.8181
$ws.Range("L122").Value = 11950.7499
$ws.Range("M122").Value = -6241.8181
$ws.Range("N122").Value = -16850.7499
$ws.Range("H132").Value = 4853.3657
$ws.Range("I132").Value = 1720
$ws.Range("J132").Value = 9276.941
$ws.Range("K132").Value = 5160
$ws.Range("L132").Value = 27830.823
$ws.Range("M132").Value = -2630
$ws.Range("N132").Value = -32890.823

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 19659
$ws.Range("J57").Value = 19659
$ws.Range("L57").Value = 19659
$ws.Range("N57").Value = -21167
$ws.Range("H80").Value = 39720.2
$ws.Range("J80").Value = 39720.2
$ws.Range("L80").Value = 39720.2
$ws.Range("N80").Value = -41716.2
$ws.Range("H83").Value = 39720.2
$ws.Range("J83").Value = 39720.2
$ws.Range("L83").Value = 119160.6
$ws.Range("N83").Value = -129144.6
$ws.Range("H123").Value = 29965
$ws.Range("J123").Value = 29965
$ws.Range("L123").Value = 29965
$ws.Range("N123").Value = -39765
